# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
# Column D ("Price") values that look like plain numbers are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# inlineStr/shared-string cells already in the sheet) instead of silently
# parsing them into floats and losing trailing zeros / exact formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.989.45"
$ws.Range("E2").Value = "  -5.33%  "

# Row 3
$ws.Range("D3").Value = "2.990.46"
$ws.Range("E3").Value = "  -5.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "'568.11"
$ws.Range("E5").Value = "  -5.73%  "

# Row 6
$ws.Range("D6").Value = "'124.22"
$ws.Range("E6").Value = "  -8.42%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").Value = "2.980.11"
$ws.Range("E8").Value = "  -5.90%  "

# Row 9
$ws.Range("D9").Value = "'0.500"
$ws.Range("E9").Value = "  -2.54%  "

# Row 10
$ws.Range("E10").Value = "  -8.36%  "

# Row 11
$ws.Range("D11").Value = "'4.95"
$ws.Range("E11").Value = "  -7.40%  "

# Row 12
$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = "  -3.00%  "

# Row 13
$ws.Range("D13").Value = "'0.0000219"
$ws.Range("E13").Value = "  -8.61%  "

# Row 14
$ws.Range("D14").Value = "'32.31"
$ws.Range("E14").Value = "  -7.23%  "

# Row 15
$ws.Range("E15").Value = "  +0.12%  "

# Row 16
$ws.Range("D16").Value = "3.485.56"
$ws.Range("E16").Value = "  -5.55%  "

# Row 17
$ws.Range("D17").Value = "2.988.64"
$ws.Range("E17").Value = "  -5.72%  "

# Row 18
$ws.Range("D18").Value = "60.010.90"
$ws.Range("E18").Value = "  -5.21%  "

# Row 19
$ws.Range("D19").Value = "'6.47"
$ws.Range("E19").Value = "  -1.51%  "

# Row 20
$ws.Range("D20").Value = "'425.63"
$ws.Range("E20").Value = "  -7.78%  "

# Row 21
$ws.Range("D21").Value = "'13.10"
$ws.Range("E21").Value = "  -6.40%  "

# Row 22
$ws.Range("D22").Value = "'0.667"
$ws.Range("E22").Value = "  -4.35%  "

# Row 23
$ws.Range("D23").Value = "'7.01"
$ws.Range("E23").Value = "  -8.53%  "

# Row 24
$ws.Range("D24").Value = "'12.79"
$ws.Range("E24").Value = "  -3.31%  "

# Row 25
$ws.Range("D25").Value = "'79.21"
$ws.Range("E25").Value = "  -4.69%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("E27").Value = "  +0.10%  "

# Row 28
$ws.Range("D28").Value = "'2.51"
$ws.Range("E28").Value = "  -7.01%  "

# Row 29
$ws.Range("D29").Value = "'1.93"
$ws.Range("E29").Value = "  -6.79%  "

# Row 30
$ws.Range("D30").Value = "'7.14"
$ws.Range("E30").Value = "  -7.56%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.07"
$ws.Range("E31").Value = "  -10.45%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'25.12"
$ws.Range("E32").Value = "  -7.57%  "

# Row 33
$ws.Range("D33").Value = "'0.0962"
$ws.Range("E33").Value = "  -4.44%  "

# Row 34
$ws.Range("D34").Value = "'5.56"
$ws.Range("E34").Value = "  -5.50%  "

# Row 35
$ws.Range("D35").Value = "'0.921"
$ws.Range("E35").Value = "  -9.44%  "

# Row 36
$ws.Range("D36").Value = "'50.13"
$ws.Range("E36").Value = "  -2.54%  "

# Row 37
$ws.Range("D37").Value = "'1.96"
$ws.Range("E37").Value = "  -19.23%  "

# Row 38
$ws.Range("E38").Value = "  +4.18%  "

# Row 39
$ws.Range("D39").Value = "0.0₃0650"
$ws.Range("E39").Value = "  -11.16%  "

# Row 40
$ws.Range("D40").Value = "'0.0352"
$ws.Range("E40").Value = "  -9.55%  "

# Row 41
$ws.Range("E41").Value = "  -5.30%  "

# Row 42
$ws.Range("D42").Value = "'368.72"
$ws.Range("E42").Value = "  -6.13%  "

# Row 43
$ws.Range("D43").Value = "2.660.25"
$ws.Range("E43").Value = "  -4.75%  "

# Row 44
$ws.Range("D44").Value = "'2.41"
$ws.Range("E44").Value = "  -8.38%  "

# Row 46
$ws.Range("D46").Value = "'120.55"
$ws.Range("E46").Value = "  -4.95%  "

# Row 47
$ws.Range("D47").Value = "'0.232"
$ws.Range("E47").Value = "  -7.28%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.107"
$ws.Range("E48").Value = "  -3.63%  "

# Row 49
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").Value = "'1.96"
$ws.Range("E49").Value = "  -7.07%  "

# Row 50
$ws.Range("D50").Value = "'23.16"
$ws.Range("E50").Value = "  -7.76%  "

# Row 51
$ws.Range("D51").Value = "'1.99"
$ws.Range("E51").Value = "  -7.94%  "
